$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interested count) in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 10388
$wsExpo.Range("F3").Value = 430
$wsExpo.Range("F4").Value = 2531
$wsExpo.Range("F6").Value = 286
$wsExpo.Range("F7").Value = 190
$wsExpo.Range("F9").Value = 778
$wsExpo.Range("F12").Value = 1092
$wsExpo.Range("F13").Value = 3239
$wsExpo.Range("F14").Value = 2410
$wsExpo.Range("F16").Value = 2161
$wsExpo.Range("F17").Value = 2161
$wsExpo.Range("F21").Value = 1595
$wsExpo.Range("F22").Value = 577
$wsExpo.Range("F23").Value = 68
$wsExpo.Range("F24").Value = 250
$wsExpo.Range("F26").Value = 24
$wsExpo.Range("F27").Value = 242
$wsExpo.Range("F29").Value = 382
$wsExpo.Range("F30").Value = 7
$wsExpo.Range("F32").Value = 394
$wsExpo.Range("F34").Value = 18
$wsExpo.Range("F36").Value = 261
$wsExpo.Range("F39").Value = 474
$wsExpo.Range("F40").Value = 459
$wsExpo.Range("F42").Value = 140
$wsExpo.Range("F46").Value = 1029
$wsExpo.Range("F48").Value = 366

# Sheet "全部类型" (All Types) - same underlying data, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10388
$wsAll.Range("F3").Value = 430
$wsAll.Range("F4").Value = 2531
$wsAll.Range("F8").Value = 286
$wsAll.Range("F9").Value = 190
$wsAll.Range("F11").Value = 778
$wsAll.Range("F12").Value = 1092
$wsAll.Range("F13").Value = 3239
$wsAll.Range("F14").Value = 2410
$wsAll.Range("F15").Value = 2161
$wsAll.Range("F16").Value = 2161
$wsAll.Range("F17").Value = 1595
$wsAll.Range("F18").Value = 577
$wsAll.Range("F19").Value = 68
$wsAll.Range("F20").Value = 250
$wsAll.Range("F22").Value = 24
$wsAll.Range("F23").Value = 242
$wsAll.Range("F25").Value = 382
$wsAll.Range("F26").Value = 7
$wsAll.Range("F28").Value = 394
$wsAll.Range("F30").Value = 18
$wsAll.Range("F35").Value = 261
$wsAll.Range("F37").Value = 474
$wsAll.Range("F39").Value = 459
$wsAll.Range("F41").Value = 140
$wsAll.Range("F48").Value = 1029
$wsAll.Range("F49").Value = 366
